# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the per-row records (columns D, H, I, J, K,
# L, M, N, P, Q) among the existing data rows 2-18; columns A, B, C, E, F,
# G, O, R are identical across every row so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "record" per row.
$cols = @(4, 8, 9, 10, 11, 12, 13, 14, 16, 17)   # D,H,I,J,K,L,M,N,P,Q

# Snapshot the current (pre-edit) values for every data row so the
# permutation below reads consistently-old data regardless of write order.
$snapshot = @{}
for ($r = 2; $r -le 18; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        # .Value2 (not .Value) — this COM shim mis-reads .Value as its own
        # property-descriptor wrapper for some cell kinds.
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# destinationRow -> sourceRow (pre-edit row whose record now lives here)
$mapping = @{
    2  = 13
    3  = 6
    4  = 14
    5  = 7
    6  = 10
    7  = 11
    8  = 17
    9  = 3
    10 = 2
    11 = 4
    12 = 15
    13 = 5
    14 = 12
    15 = 16
    16 = 18
    17 = 8
    18 = 9
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
